# edit.ps1 - applies the curso_1.docx form-field edits described by the diff.
#
# The document body is dominated by one large form table (the 3rd table in
# the document) whose rows are label/value pairs. Most of the diff hunks
# are simple value replacements in that table; a handful of hunks move an
# "X" checkbox mark between cells in the "Dias semana de programacion" row,
# and the final two hunks touch the instructor signature paragraph near the
# end of the document.

$d = $word.ActiveDocument
$t = $d.Tables.Item(3)

function Set-CellText($table, $rowIndex, $colIndex, $newText) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item($colIndex)
    $cellRange = $cell.Range
    $target = $d.Range($cellRange.Start, $cellRange.End - 1)
    $target.Text = $newText
}

function Clear-CellText($table, $rowIndex, $colIndex) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item($colIndex)
    $cellRange = $cell.Range
    $target = $d.Range($cellRange.Start, $cellRange.End - 1)
    $target.Text = ""
}

# --- Simple field value replacements -------------------------------------

# Codigo programa de formacion*: 12 -> 8220041
Set-CellText $t 1 2 "8220041"

# Nombre del Programa*: sistemas -> TRIBOLOGIA Y LUBRICACION
Set-CellText $t 2 2 "TRIBOLOGIA Y LUBRICACION"

# Version del programa*: 1 -> 08
Set-CellText $t 3 2 "08"

# Fecha de Inicio*: 2025-09-22 -> 2025-09-01
Set-CellText $t 5 2 "2025-09-01"

# Fecha prevista de terminacion*: 2025-09-30 -> 2025-10-30
Set-CellText $t 6 2 "2025-10-30"

# Departamento desarrollo de formacion*: Caldas -> Cauca
Set-CellText $t 9 2 "Cauca"

# Municipio desarrollo formacion*: La Dorada -> Popayan
Set-CellText $t 10 2 "Popayán"

# Direccion donde se va a realizar la formacion*: tics1 -> asd123
Set-CellText $t 11 2 "asd123"

# Nombre responsable*: GIlber Martinez -> Instructor Sena
Set-CellText $t 12 2 "Instructor Sena"

# (tipo de documento): SIN DOCUMENTO # -> CC #
Set-CellText $t 12 3 "CC #"

# (numero de documento): j6mF77U0 -> 234234234
Set-CellText $t 12 4 "234234234"

# Correo electronico*: 1@gmail.com -> instructor@gmail.com
Set-CellText $t 13 2 "instructor@gmail.com"

# --- "Dias semana de programacion*" checkbox row --------------------------
# Remove the X marks for MAR (col 4) and JUE (col 8); add an X for VIE (col 10)
Clear-CellText $t 33 4
Clear-CellText $t 33 8
Set-CellText $t 33 10 "X"

# --- Horario / fechas de ejecucion ----------------------------------------

# Horario del curso de formacion*: "7:00 15:00" -> "8 a 12" (only the first
# run of this cell; the remainder of the cell text is left untouched)
$horarioCell = $t.Rows.Item(34).Cells.Item(2)
$oldHorario = "7:00 15:00"
$horarioRange = $d.Range($horarioCell.Range.Start, $horarioCell.Range.Start + $oldHorario.Length)
$horarioRange.Text = "8 a 12"

# Fechas de ejecucion de la formacion (mes 1): "14 15 16 17" -> "12"
Set-CellText $t 35 2 "12"

# Fechas de ejecucion de la formacion (mes 2): "1" -> "22"
Set-CellText $t 36 2 "22"

# --- Footer signature paragraph -------------------------------------------
# "Nombre del instructor:GIlber Martinez <tabs>Firma Instructor sena"
$footerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Nombre del instructor*") {
        $footerPara = $p
        break
    }
}

$words = $footerPara.Range.Words
$nameWord = $words.Item(5)
$oldName = "GIlber Martinez "
$nameRange = $d.Range($nameWord.Start, $nameWord.Start + $oldName.Length)
$nameRange.Text = "Instructor Sena "

# Re-fetch the paragraph/words since the document shifted after the edit above
$footerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Firma Instructor*") {
        $footerPara = $p
        break
    }
}
$words2 = $footerPara.Range.Words
$lastWord = $words2.Item($words2.Count - 1)
$lastWord.Text = "qwe"

Write-Output "done"
